$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = -17.89928610827413
$ws.Cells.Item(2, 3).Value = 1.986167488036832
$ws.Cells.Item(2, 4).Value = -17.89928610827413
$ws.Cells.Item(2, 5).Value = -17.89928610827413
$ws.Cells.Item(2, 6).Value = -17.89928610827413
$ws.Cells.Item(2, 7).Value = -17.89928610827413
$ws.Cells.Item(2, 8).Value = -17.89928610827413
$ws.Cells.Item(2, 9).Value = -17.89928610827413
$ws.Cells.Item(2, 10).Value = -17.89928610827413
$ws.Cells.Item(2, 11).Value = -17.89928610827413

# Row 3
$ws.Cells.Item(3, 2).Value = -17.89928610827413
$ws.Cells.Item(3, 3).Value = -17.89928610827413
$ws.Cells.Item(3, 4).Value = -17.89928610827413
$ws.Cells.Item(3, 5).Value = -17.89928610827413
$ws.Cells.Item(3, 6).Value = -17.89928610827413
$ws.Cells.Item(3, 7).Value = -17.89928610827413
$ws.Cells.Item(3, 8).Value = -17.89928610827413
$ws.Cells.Item(3, 9).Value = 2.874143701379538
$ws.Cells.Item(3, 10).Value = -17.89928610827413
$ws.Cells.Item(3, 11).Value = -17.89928610827413

# Row 4
$ws.Cells.Item(4, 2).Value = -17.89928610827413
$ws.Cells.Item(4, 3).Value = 2.246134869814681
$ws.Cells.Item(4, 4).Value = 2.235372082233521
$ws.Cells.Item(4, 5).Value = -17.89928610827413
$ws.Cells.Item(4, 6).Value = 3.388190123347006
$ws.Cells.Item(4, 7).Value = -17.89928610827413
$ws.Cells.Item(4, 8).Value = 1.674754846788818
$ws.Cells.Item(4, 9).Value = -17.89928610827413
$ws.Cells.Item(4, 10).Value = 2.531289970499262
$ws.Cells.Item(4, 11).Value = -17.89928610827413

# Row 5
$ws.Cells.Item(5, 2).Value = -17.89928610827413
$ws.Cells.Item(5, 3).Value = 1.850742034598682
$ws.Cells.Item(5, 4).Value = -17.89928610827413
$ws.Cells.Item(5, 5).Value = -17.89928610827413
$ws.Cells.Item(5, 6).Value = -17.89928610827413
$ws.Cells.Item(5, 7).Value = -17.89928610827413
$ws.Cells.Item(5, 8).Value = -17.89928610827413
$ws.Cells.Item(5, 9).Value = -17.89928610827413
$ws.Cells.Item(5, 10).Value = -17.89928610827413
$ws.Cells.Item(5, 11).Value = -17.89928610827413

# Row 6
$ws.Cells.Item(6, 2).Value = -17.89928610827413
$ws.Cells.Item(6, 3).Value = -17.89928610827413
$ws.Cells.Item(6, 4).Value = -17.89928610827413
$ws.Cells.Item(6, 5).Value = -17.89928610827413
$ws.Cells.Item(6, 6).Value = -17.89928610827413
$ws.Cells.Item(6, 7).Value = -17.89928610827413
$ws.Cells.Item(6, 8).Value = -17.89928610827413
$ws.Cells.Item(6, 9).Value = -17.89928610827413
$ws.Cells.Item(6, 10).Value = -17.89928610827413
$ws.Cells.Item(6, 11).Value = -17.89928610827413

# Row 7
$ws.Cells.Item(7, 2).Value = 2.628188034202372
$ws.Cells.Item(7, 3).Value = -17.89928610827413
$ws.Cells.Item(7, 4).Value = -17.89928610827413
$ws.Cells.Item(7, 5).Value = -17.89928610827413
$ws.Cells.Item(7, 6).Value = -17.89928610827413
$ws.Cells.Item(7, 7).Value = -17.89928610827413
$ws.Cells.Item(7, 8).Value = -17.89928610827413
$ws.Cells.Item(7, 9).Value = -17.89928610827413
$ws.Cells.Item(7, 10).Value = -17.89928610827413
$ws.Cells.Item(7, 11).Value = -17.89928610827413

# Row 8
$ws.Cells.Item(8, 2).Value = -17.89928610827413
$ws.Cells.Item(8, 3).Value = -17.89928610827413
$ws.Cells.Item(8, 4).Value = -17.89928610827413
$ws.Cells.Item(8, 5).Value = 1.813719396038907
$ws.Cells.Item(8, 6).Value = -17.89928610827413
$ws.Cells.Item(8, 7).Value = -17.89928610827413
$ws.Cells.Item(8, 8).Value = -17.89928610827413
$ws.Cells.Item(8, 9).Value = -17.89928610827413
$ws.Cells.Item(8, 10).Value = -17.89928610827413
$ws.Cells.Item(8, 11).Value = -17.89928610827413

# Row 9
$ws.Cells.Item(9, 2).Value = 3.788418049035114
$ws.Cells.Item(9, 3).Value = -17.89928610827413
$ws.Cells.Item(9, 4).Value = -17.89928610827413
$ws.Cells.Item(9, 5).Value = -17.89928610827413
$ws.Cells.Item(9, 6).Value = -17.89928610827413
$ws.Cells.Item(9, 7).Value = -17.89928610827413
$ws.Cells.Item(9, 8).Value = -17.89928610827413
$ws.Cells.Item(9, 9).Value = -17.89928610827413
$ws.Cells.Item(9, 10).Value = -17.89928610827413
$ws.Cells.Item(9, 11).Value = -17.89928610827413

# Row 10
$ws.Cells.Item(10, 2).Value = -17.89928610827413
$ws.Cells.Item(10, 3).Value = -17.89928610827413
$ws.Cells.Item(10, 4).Value = -17.89928610827413
$ws.Cells.Item(10, 5).Value = -17.89928610827413
$ws.Cells.Item(10, 6).Value = -17.89928610827413
$ws.Cells.Item(10, 7).Value = -17.89928610827413
$ws.Cells.Item(10, 8).Value = -17.89928610827413
$ws.Cells.Item(10, 9).Value = 1.4377542677831
$ws.Cells.Item(10, 10).Value = -17.89928610827413
$ws.Cells.Item(10, 11).Value = 1.867658449561808

# Row 11
$ws.Cells.Item(11, 2).Value = -17.89928610827413
$ws.Cells.Item(11, 3).Value = -17.89928610827413
$ws.Cells.Item(11, 4).Value = -17.89928610827413
$ws.Cells.Item(11, 5).Value = 2.861759428749952
$ws.Cells.Item(11, 6).Value = -17.89928610827413
$ws.Cells.Item(11, 7).Value = -17.89928610827413
$ws.Cells.Item(11, 8).Value = -17.89928610827413
$ws.Cells.Item(11, 9).Value = -17.89928610827413
$ws.Cells.Item(11, 10).Value = -17.89928610827413
$ws.Cells.Item(11, 11).Value = 1.693267649507455

# Row 12
$ws.Cells.Item(12, 2).Value = -17.89928610827413
$ws.Cells.Item(12, 3).Value = -17.89928610827413
$ws.Cells.Item(12, 4).Value = -17.89928610827413
$ws.Cells.Item(12, 5).Value = -17.89928610827413
$ws.Cells.Item(12, 6).Value = -17.89928610827413
$ws.Cells.Item(12, 7).Value = -17.89928610827413
$ws.Cells.Item(12, 8).Value = -17.89928610827413
$ws.Cells.Item(12, 9).Value = -17.89928610827413
$ws.Cells.Item(12, 10).Value = -17.89928610827413
$ws.Cells.Item(12, 11).Value = -17.89928610827413

# Row 13
$ws.Cells.Item(13, 2).Value = -17.89928610827413
$ws.Cells.Item(13, 3).Value = -17.89928610827413
$ws.Cells.Item(13, 4).Value = -17.89928610827413
$ws.Cells.Item(13, 5).Value = 2.361877978286752
$ws.Cells.Item(13, 6).Value = -17.89928610827413
$ws.Cells.Item(13, 7).Value = -17.89928610827413
$ws.Cells.Item(13, 8).Value = -17.89928610827413
$ws.Cells.Item(13, 9).Value = -17.89928610827413
$ws.Cells.Item(13, 10).Value = 2.207865822622836
$ws.Cells.Item(13, 11).Value = 2.027022918243478

# Row 14
$ws.Cells.Item(14, 2).Value = -17.89928610827413
$ws.Cells.Item(14, 3).Value = -17.89928610827413
$ws.Cells.Item(14, 4).Value = 1.283149468279633
$ws.Cells.Item(14, 5).Value = -17.89928610827413
$ws.Cells.Item(14, 6).Value = -17.89928610827413
$ws.Cells.Item(14, 7).Value = -17.89928610827413
$ws.Cells.Item(14, 8).Value = -17.89928610827413
$ws.Cells.Item(14, 9).Value = -17.89928610827413
$ws.Cells.Item(14, 10).Value = -17.89928610827413
$ws.Cells.Item(14, 11).Value = 2.170731296051785

# Row 15
$ws.Cells.Item(15, 2).Value = -17.89928610827413
$ws.Cells.Item(15, 3).Value = -17.89928610827413
$ws.Cells.Item(15, 4).Value = 1.246956878257685
$ws.Cells.Item(15, 5).Value = -17.89928610827413
$ws.Cells.Item(15, 6).Value = -17.89928610827413
$ws.Cells.Item(15, 7).Value = -17.89928610827413
$ws.Cells.Item(15, 8).Value = -17.89928610827413
$ws.Cells.Item(15, 9).Value = -17.89928610827413
$ws.Cells.Item(15, 10).Value = -17.89928610827413
$ws.Cells.Item(15, 11).Value = -17.89928610827413

# Row 16
$ws.Cells.Item(16, 2).Value = -17.89928610827413
$ws.Cells.Item(16, 3).Value = -17.89928610827413
$ws.Cells.Item(16, 4).Value = -17.89928610827413
$ws.Cells.Item(16, 5).Value = -17.89928610827413
$ws.Cells.Item(16, 6).Value = -17.89928610827413
$ws.Cells.Item(16, 7).Value = -17.89928610827413
$ws.Cells.Item(16, 8).Value = -17.89928610827413
$ws.Cells.Item(16, 9).Value = -17.89928610827413
$ws.Cells.Item(16, 10).Value = 2.233892486191674
$ws.Cells.Item(16, 11).Value = -17.89928610827413

# Row 17
$ws.Cells.Item(17, 2).Value = -17.89928610827413
$ws.Cells.Item(17, 3).Value = 1.882971487983077
$ws.Cells.Item(17, 4).Value = 2.151839173704065
$ws.Cells.Item(17, 5).Value = -17.89928610827413
$ws.Cells.Item(17, 6).Value = -17.89928610827413
$ws.Cells.Item(17, 7).Value = -17.89928610827413
$ws.Cells.Item(17, 8).Value = 1.052172216033344
$ws.Cells.Item(17, 9).Value = 1.435225063111893
$ws.Cells.Item(17, 10).Value = 1.181440352996278
$ws.Cells.Item(17, 11).Value = -17.89928610827413

# Row 18
$ws.Cells.Item(18, 2).Value = -17.89928610827413
$ws.Cells.Item(18, 3).Value = -17.89928610827413
$ws.Cells.Item(18, 4).Value = -17.89928610827413
$ws.Cells.Item(18, 5).Value = -17.89928610827413
$ws.Cells.Item(18, 6).Value = -17.89928610827413
$ws.Cells.Item(18, 7).Value = -17.89928610827413
$ws.Cells.Item(18, 8).Value = 0.9319252479547203
$ws.Cells.Item(18, 9).Value = 0.8431976797712755
$ws.Cells.Item(18, 10).Value = 1.393469134322878
$ws.Cells.Item(18, 11).Value = -17.89928610827413

# Row 19
$ws.Cells.Item(19, 2).Value = -17.89928610827413
$ws.Cells.Item(19, 3).Value = -17.89928610827413
$ws.Cells.Item(19, 4).Value = 1.65468997449406
$ws.Cells.Item(19, 5).Value = -17.89928610827413
$ws.Cells.Item(19, 6).Value = -17.89928610827413
$ws.Cells.Item(19, 7).Value = -17.89928610827413
$ws.Cells.Item(19, 8).Value = 1.704139212568323
$ws.Cells.Item(19, 9).Value = 1.775411337496887
$ws.Cells.Item(19, 10).Value = -17.89928610827413
$ws.Cells.Item(19, 11).Value = -17.89928610827413

# Row 20
$ws.Cells.Item(20, 2).Value = -17.89928610827413
$ws.Cells.Item(20, 3).Value = 0.8928591646948326
$ws.Cells.Item(20, 4).Value = 1.531821730900515
$ws.Cells.Item(20, 5).Value = -17.89928610827413
$ws.Cells.Item(20, 6).Value = 3.252464423740694
$ws.Cells.Item(20, 7).Value = -17.89928610827413
$ws.Cells.Item(20, 8).Value = 2.069680146692975
$ws.Cells.Item(20, 9).Value = 1.026908524537942
$ws.Cells.Item(20, 10).Value = -17.89928610827413
$ws.Cells.Item(20, 11).Value = 2.182219434857289

# Row 21
$ws.Cells.Item(21, 2).Value = -17.89928610827413
$ws.Cells.Item(21, 3).Value = 1.098932427886121
$ws.Cells.Item(21, 4).Value = -17.89928610827413
$ws.Cells.Item(21, 5).Value = 2.026820189014712
$ws.Cells.Item(21, 6).Value = -17.89928610827413
$ws.Cells.Item(21, 7).Value = 4.321922488579167
$ws.Cells.Item(21, 8).Value = 2.424798824627789
$ws.Cells.Item(21, 9).Value = -17.89928610827413
$ws.Cells.Item(21, 10).Value = -17.89928610827413
$ws.Cells.Item(21, 11).Value = -17.89928610827413

